# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell (A1) onto the
# three new header cells so they match the bold/bordered/centered style
# used by every other column header.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-50) gets the same season record values.
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 87   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 75   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
